$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the planning date shown in C1 --------------------------------
# The cell is styled with a date number format (numFmtId 14) but actually
# stores its value as literal text ("28/09/2020"), so a direct
# Range("C1").Value = "01/10/2020" would be auto-coerced by Excel into a
# real date serial number. To preserve the original "text that looks like
# a date" semantics (and keep the existing cell style untouched), we stage
# the text in a scratch cell explicitly formatted as Text, copy/paste the
# value into C1, then wipe the scratch cell completely.
$scratch = $ws.Range("Z1")
$scratch.NumberFormat = "@"
$scratch.Value = "01/10/2020"
$scratch.Copy()
$ws.Range("C1").PasteSpecial(-4163)
$scratch.Clear()

# --- Row 3 (Chambre A / Ethan) --------------------------------------------
$ws.Range("B3").Value = "Ethan"
$ws.Range("C3").ClearContents()
$ws.Range("E3").Value = "Anne BOROCCO`nNo desc"
$ws.Range("G3").Value = "Danièle LEVY`nNo desc"

# --- Row 4 (Chambre B / Suzon) --------------------------------------------
$ws.Range("C4").ClearContents()
$ws.Range("D4").Value = "Paula HARRY`nNo desc"
$ws.Range("E4").Value = "Sylvie COHEN`nNo desc"
$ws.Range("F4").Value = "Danièle LEVY`nNo desc"

# --- Row 5 (Chambre C / Maïmouna) -----------------------------------------
$ws.Range("B5").Value = "Maïmouna"
$ws.Range("D5").ClearContents()
$ws.Range("E5").Value = "Claude ARDITTY`nNo desc"
$ws.Range("H5").Value = "Sylvie COHEN`nNo desc"

# --- Row 6 (Chambre D) ------------------------------------------------------
$ws.Range("B6").ClearContents()
$ws.Range("E6").ClearContents()

# --- Row 7 (Chambre E / Maxence) -------------------------------------------
$ws.Range("C7").Value = "Claude ARDITTY`nNo desc"
$ws.Range("D7").Value = "Anne BOROCCO`nNo desc"
$ws.Range("E7").Value = "Danièle LEVY`nNo desc"
$ws.Range("F7").Value = "Paula HARRY`nNo desc"

# --- Row 8 (Chambre F / Salim) ----------------------------------------------
$ws.Range("C8").Value = "Stéphanie DESQUIENS`nNo desc"
$ws.Range("H8").Value = "Anne BOROCCO`nNo desc"

# --- Row 9 (Chambre G / Mégane) ---------------------------------------------
$ws.Range("B9").Value = "Mégane"
$ws.Range("D9").Value = "Claude ARDITTY`nNo desc"
$ws.Range("F9").Value = "Sylvie COHEN`nNo desc"

# --- Row 10 (Chambre H / Nasreddine) -----------------------------------------
$ws.Range("B10").Value = "Nasreddine"

# --- Row 11 (Chambre I / Ninon) -----------------------------------------------
$ws.Range("D11").Value = "Stéphanie DESQUIENS`nNo desc"
$ws.Range("E11").Value = "Paula HARRY`nNo desc"
$ws.Range("G11").Value = "Sylvie COHEN`nNo desc"

# --- Row 12 (Chambre J / Carmen) ----------------------------------------------
$ws.Range("D12").Value = "Sylvie COHEN`nNo desc"
$ws.Range("E12").Value = "Stéphanie DESQUIENS`nNo desc"
$ws.Range("H12").Value = "Paula HARRY`nNo desc"

# --- Row 13 (Chambre K / Charlène) --------------------------------------------
$ws.Range("F13").Value = "Claude ARDITTY`nNo desc"
$ws.Range("G13").Value = "Paula HARRY`nNo desc"
$ws.Range("H13").Value = "Danièle LEVY`nNo desc"

# --- Row 14 (Chambre L) ---------------------------------------------------------
$ws.Range("B14").ClearContents()

# --- Row 15 (HDJ 1 / Jason) ------------------------------------------------------
$ws.Range("A15").Value = "HDJ 1"
$ws.Range("B15").Value = "Jason"
$ws.Range("E15").ClearContents()
$ws.Range("F15").Value = "Anne BOROCCO`nNo desc"

# --- Row 16 (HDJ 2) ---------------------------------------------------------------
$ws.Range("A16").Value = "HDJ 2"
$ws.Range("B16").ClearContents()

# --- Row 17 (HDJ 3) ---------------------------------------------------------------
$ws.Range("A17").Value = "HDJ 3"
$ws.Range("B17").ClearContents()
